# Update "想去人数" (interested-count) figures that changed between the
# previous and newly generated gh-pages data snapshot.
#
# Sheet "展览"   (展览/exhibition)
# Sheet "演出"   (演出/show)
# Sheet "本地生活" (local life)
# Sheet "全部类型" (all types - combined listing)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 133
$ws1.Range("F11").Value = 5685
$ws1.Range("F12").Value = 40
$ws1.Range("F16").Value = 539
$ws1.Range("F17").Value = 334
$ws1.Range("F25").Value = 61
$ws1.Range("F26").Value = 1739

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 46

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 189

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 189
$ws4.Range("F10").Value = 133
$ws4.Range("F13").Value = 5685
$ws4.Range("F14").Value = 40
$ws4.Range("F19").Value = 539
$ws4.Range("F20").Value = 334
$ws4.Range("F22").Value = 46
$ws4.Range("F35").Value = 61
$ws4.Range("F36").Value = 1739
